$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the building name labels (shared strings) in column A for existing rows
$ws.Range("A2").Value = "building/ghx_GHX_shkwfn"
$ws.Range("A3").Value = "building/ghx_GHP_building_mksuwer"
$ws.Range("A4").Value = "building/ghx_GHP_building_sdupkgra"

# Update row 2 values
$ws.Range("B2").Value = 77335.30809999999
$ws.Range("C2").Value = 81472.3
$ws.Range("F2").Value = 81472.3
$ws.Range("I2").Value = 1.812
$ws.Range("J2").Value = 81472.3

# Update row 3 values
$ws.Range("B3").Value = 18062.1037
$ws.Range("C3").Value = 8320.525
$ws.Range("F3").Value = 8320.525
$ws.Range("I3").Value = 2.785
$ws.Range("J3").Value = 8320.525

# Update row 4 values
$ws.Range("B4").Value = 238067.3911
$ws.Range("C4").Value = 47962.5012
$ws.Range("F4").Value = 47962.5012
$ws.Range("I4").Value = 1.384
$ws.Range("J4").Value = 47962.5012

# Add new row 5
$ws.Range("A5").Value = "building/ghx_GHX_ahudfd"
$ws.Range("B5").Value = 77335.30809999999
$ws.Range("C5").Value = 81472.3
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 42.15
$ws.Range("F5").Value = 81472.3
$ws.Range("G5").Value = 8261.8833
$ws.Range("H5").Value = 66.18000000000001
$ws.Range("I5").Value = 1.398
$ws.Range("J5").Value = 81472.3
$ws.Range("K5").Value = 437.6631368059719
$ws.Range("L5").Value = 761.27
$ws.Range("M5").Value = 1.68
$ws.Range("N5").Value = 0.0672
